$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 205: Total row - sums of the True Positive / False Positive / False Negative columns
$ws.Range("A205").Value = "Total "
$ws.Range("B205").Formula = "=SUM(B2:B203)"
$ws.Range("C205").Formula = "=SUM(C2:C203)"
$ws.Range("D205").Formula = "=SUM(D2:D203)"

# Row 206: Recall
$ws.Range("A206").Value = "Recall"
$ws.Range("B206").Formula = "=B205/(16*202)"

# Row 207: Precision
$ws.Range("A207").Value = "Precision"
$ws.Range("B207").Formula = "=B205/(C205+B205)"

# Row 208: F1 score
$ws.Range("A208").Value = "F1"
$ws.Range("B208").Formula = "=2*B206*B207/(B206+B207)"

# Row 209: label the existing "Images with perfect score" COUNTIF cell (B209 already holds
# the formula / value, we are only adding the caption in column A)
$ws.Range("A209").Value = "Images with perfect score"

# Widen column A slightly so the new labels fit
$ws.Columns.Item(1).ColumnWidth = 21.83

# Scroll the view down to the newly added summary rows and move the selection there
$ws.Range("A195").Select()
$excel.ActiveWindow.DisplayGridlines = $true
